$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.418.22"
$ws.Range("E2").Value = "  +2.01%  "
$ws.Range("D3").Value = "2.664.76"
$ws.Range("E3").Value = "  +2.04%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'602.51"
$ws.Range("E5").Value = "  -0.42%  "
$ws.Range("D6").Value = "'178.84"
$ws.Range("E6").Value = "  -0.55%  "
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").Value = "'0.524"
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("E9").Value = "  +6.43%  "
$ws.Range("D10").Value = "2.662.70"
$ws.Range("E10").Value = "  +2.02%  "
$ws.Range("D11").Value = "'0.168"
$ws.Range("D12").Value = "'0.358"
$ws.Range("E12").Value = "  +3.48%  "
$ws.Range("D13").Value = "'5.03"
$ws.Range("E13").Value = "  -0.13%  "
$ws.Range("D14").Value = "3.148.33"
$ws.Range("E14").Value = "  +1.59%  "
$ws.Range("D15").Value = "'0.0000190"
$ws.Range("E15").Value = "  +4.17%  "
$ws.Range("D16").Value = "72.287.86"
$ws.Range("E16").Value = "  +1.79%  "
$ws.Range("D17").Value = "'26.58"
$ws.Range("E17").Value = "  -0.07%  "
$ws.Range("D18").Value = "2.661.31"
$ws.Range("E18").Value = "  +1.65%  "
$ws.Range("D19").Value = "'12.00"
$ws.Range("E19").Value = "  +4.56%  "
$ws.Range("D20").Value = "'8.02"
$ws.Range("E20").Value = "  +2.25%  "
$ws.Range("D21").Value = "'379.17"
$ws.Range("E21").Value = "  +0.18%  "
$ws.Range("E22").Value = "  +2.40%  "
$ws.Range("E23").Value = "  +12.38%  "
$ws.Range("D24").Value = "'72.70"
$ws.Range("E24").Value = "  +1.12%  "
$ws.Range("E25").Value = "  +0.52%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").Value = "'9.99"
$ws.Range("E27").Value = "  +4.72%  "
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("D30").Value = "0.0₃0955"
$ws.Range("E30").Value = "  +0.56%  "
$ws.Range("D31").Value = "'8.25"
$ws.Range("E31").Value = "  +3.27%  "
$ws.Range("D32").Value = "'524.59"
$ws.Range("E32").Value = "  -1.00%  "
$ws.Range("E33").Value = "  +0.08%  "
$ws.Range("E34").Value = "  -0.16%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  -0.18%  "
$ws.Range("D36").Value = "'163.34"
$ws.Range("E36").Value = "  -1.29%  "
$ws.Range("D37").Value = "'19.53"
$ws.Range("E37").Value = "  +2.01%  "
$ws.Range("E38").Value = "  +0.83%  "
$ws.Range("E39").Value = "  -6.00%  "
$ws.Range("D40").Value = "'1.40"
$ws.Range("E40").Value = "  +1.82%  "
$ws.Range("E41").Value = "  -1.11%  "
$ws.Range("D42").Value = "'5.08"
$ws.Range("E42").Value = "  +1.32%  "
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("D44").Value = "'2.61"
$ws.Range("E44").Value = "  +0.93%  "
$ws.Range("D45").Value = "'0.336"
$ws.Range("E45").Value = "  +0.86%  "
$ws.Range("D46").Value = "'39.30"
$ws.Range("E46").Value = "  -2.01%  "
$ws.Range("D47").Value = "'153.36"
$ws.Range("E47").Value = "  -0.18%  "
$ws.Range("D48").Value = "'3.75"
$ws.Range("E48").Value = "  +2.76%  "
$ws.Range("D49").Value = "'0.549"
$ws.Range("E49").Value = "  +3.65%  "
$ws.Range("E50").Value = "  +3.00%  "
